$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # Exact, case-sensitive, whole-text match/replace anywhere in the document body.
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Header date line.
Replace-Text "2024-10-31 Thursday" "2024-11-01 Friday"

# Table of "three digit divided by one digit" problems (5 columns x 5 data rows).
# NOTE ON ORDER: one new value ("979÷7=", produced from "173÷3=") collides with an
# *old* value that appears earlier in the document ("979÷7=" -> "605÷7="). To avoid
# a just-written replacement being clobbered by a later rule, every mapping whose
# "old" text equals another mapping's "new" text is applied first.
Replace-Text "296÷2=" "405÷4="
Replace-Text "906÷5=" "353÷4="
Replace-Text "400÷9=" "772÷7="
Replace-Text "214÷4=" "457÷5="
Replace-Text "129÷6=" "357÷4="

Replace-Text "905÷3=" "268÷9="
Replace-Text "797÷6=" "243÷6="
Replace-Text "315÷3=" "851÷8="
Replace-Text "352÷6=" "249÷6="
Replace-Text "979÷7=" "605÷7="  # must run before "173÷3=" -> "979÷7=" below

Replace-Text "133÷6=" "765÷2="
Replace-Text "142÷6=" "635÷6="
Replace-Text "222÷9=" "876÷3="
Replace-Text "263÷6=" "871÷6="
Replace-Text "966÷6=" "802÷7="

Replace-Text "471÷8=" "565÷2="
Replace-Text "859÷6=" "818÷5="
Replace-Text "885÷3=" "773÷5="
Replace-Text "678÷2=" "694÷7="
Replace-Text "447÷2=" "191÷9="

Replace-Text "938÷9=" "972÷3="
Replace-Text "218÷5=" "940÷5="
Replace-Text "152÷5=" "341÷6="
Replace-Text "173÷3=" "979÷7="
Replace-Text "975÷2=" "879÷8="

Write-Output "Done"
